# Applies the "Updated cryptos list" price/volume refresh to sheet1.
# Generated from the authoritative OOXML diff: for each changed cell we
# write the new literal text. Values that parse as a pure number (e.g.
# "0.9973", "0.00001099") are written with a leading apostrophe so the
# engine keeps them as text (matching the source t="inlineStr" cells)
# instead of silently converting them to floats/doubles; the style is
# then reset to "Normal" so no stray numeric format sticks to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.436.51'
$ws.Range('E2').Value = '  -0.55%  '

$ws.Range('D3').Value = '1.790.20'
$ws.Range('E3').Value = '  -1.49%  '

$ws.Range('E4').Value = '  -0.23%  '

$ws.Range('D5').Value = "'338.80"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.13%  '

$ws.Range('D6').Value = "'0.9973"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.47%  '

$ws.Range('D7').Value = "'0.3910"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.40%  '

$ws.Range('D8').Value = "'0.3462"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.59%  '

$ws.Range('D9').Value = "'48.43"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.02%  '

$ws.Range('E10').Value = '  -2.78%  '

$ws.Range('D11').Value = "'0.07489"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.74%  '

$ws.Range('D12').Value = "'0.9964"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.73%  '

$ws.Range('D13').Value = "'21.85"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.02%  '

$ws.Range('D14').Value = "'6.507"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.47%  '

$ws.Range('D15').Value = '1.795.78'
$ws.Range('E15').Value = '  -1.19%  '

$ws.Range('D16').Value = "'7.149"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.08%  '

$ws.Range('D17').Value = "'0.00001099"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.63%  '

$ws.Range('D18').Value = "'0.06675"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.62%  '

$ws.Range('D19').Value = "'84.87"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.86%  '

$ws.Range('D20').Value = "'0.9975"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.44%  '

$ws.Range('D21').Value = "'17.69"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.75%  '

$ws.Range('D22').Value = "'6.569"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.72%  '

$ws.Range('D23').Value = '27.488.78'

$ws.Range('E24').Value = '  -4.96%  '

$ws.Range('D25').Value = "'2.408"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.58%  '

$ws.Range('D26').Value = "'2.515"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.58%  '

$ws.Range('D27').Value = "'21.27"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.33%  '

$ws.Range('D28').Value = "'1.464"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.50%  '

$ws.Range('D29').Value = "'156.79"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.61%  '

$ws.Range('D30').Value = '2.000.43'
$ws.Range('E30').Value = '  -0.95%  '

$ws.Range('D31').Value = "'135.26"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.28%  '

$ws.Range('D32').Value = "'4.036"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.20%  '

$ws.Range('D33').Value = "'6.041"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.31%  '

$ws.Range('D34').Value = "'0.08741"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.29%  '

$ws.Range('D35').Value = "'13.04"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.76%  '

$ws.Range('D36').Value = "'1.621"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.75%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = "'0.02414"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.43%  '

$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = "'5.453"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.35%  '

$ws.Range('D39').Value = "'0.06486"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.19%  '

$ws.Range('D40').Value = "'0.6833"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.77%  '

$ws.Range('E41').Value = '  -1.53%  '

$ws.Range('D42').Value = "'1.257"
$ws.Range('D42').Style = 'Normal'

$ws.Range('D43').Value = "'8.389"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.49%  '

$ws.Range('D44').Value = "'14.52"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.70%  '

$ws.Range('D45').Value = "'0.6387"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.29%  '

$ws.Range('D46').Value = "'0.9961"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.58%  '

$ws.Range('D47').Value = "'3.873"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.08%  '

$ws.Range('E48').Value = '  -1.31%  '

$ws.Range('D49').Value = "'131.78"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.62%  '

$ws.Range('D50').Value = "'0.07184"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.60%  '

$ws.Range('D51').Value = "'79.57"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.21%  '
